$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5:E32").ClearContents()
$ws.Range("F11:F32").ClearContents()
$ws.Range("G17:G32").ClearContents()

$ws.Range("I10").Select()
